# Refresh cryptocurrency price/volume data (and a few reordered rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'88.825.18"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "'  -1.95%  "
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = "'3.050.86"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'  -3.46%  "
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = "'  -0.32%  "
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = "'211.01"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "'  -1.39%  "
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = "'611.26"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "'  -2.61%  "
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.Value = "'0.362"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = "'  -9.49%  "
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.Value = "'0.881"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "'  +21.33%  "
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = "'  -0.11%  "
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.Value = "'3.049.11"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'  -3.26%  "
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'  +17.36%  "
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.Value = "'0.187"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = "'  +2.66%  "
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = "'0.0000237"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "'  -7.30%  "
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'5.32"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'  +0.74%  "
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.Value = "'89.358.47"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'  -1.16%  "
$cell.Style = "Normal"
$cell = $ws.Range("B16")
$cell.Value = "'WrappedliquidstakedEther2.0"
$cell.Style = "Normal"
$cell = $ws.Range("C16")
$cell.Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.Value = "'3.637.91"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "'  -3.25%  "
$cell.Style = "Normal"
$cell = $ws.Range("B17")
$cell.Value = "'Avalanche"
$cell.Style = "Normal"
$cell = $ws.Range("C17")
$cell.Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.Value = "'31.94"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'  -0.79%  "
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.Value = "'3.127.32"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'  -1.30%  "
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.Value = "'3.35"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'  +1.77%  "
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.Value = "'0.0000209"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'  -3.86%  "
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.Value = "'13.37"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "'  +0.73%  "
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.Value = "'423.66"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = "'  -0.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.Value = "'4.95"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "'  +0.63%  "
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = "'8.10"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "'  -3.83%  "
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.Value = "'5.39"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'  +2.45%  "
$cell.Style = "Normal"
$cell = $ws.Range("B26")
$cell.Value = "'Litecoin"
$cell.Style = "Normal"
$cell = $ws.Range("C26")
$cell.Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.Value = "'84.06"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "'  +4.90%  "
$cell.Style = "Normal"
$cell = $ws.Range("B27")
$cell.Value = "'Aptos"
$cell.Style = "Normal"
$cell = $ws.Range("C27")
$cell.Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.Value = "'11.83"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'  +2.24%  "
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.Value = "'3.241.05"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.Value = "'  -3.31%  "
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = "'  +0.07%  "
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = "'  +9.75%  "
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.Value = "'0.161"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = "'  +2.84%  "
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.Value = "'8.15"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = "'  -1.67%  "
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.Value = "'501.08"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = "'  -1.93%  "
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.Value = "'3.59"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.Value = "'  -10.59%  "
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.Value = "'6.56"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = "'  -4.39%  "
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.Value = "'22.50"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = "'  +1.50%  "
$cell.Style = "Normal"
$cell = $ws.Range("B37")
$cell.Value = "'PancakeSwap"
$cell.Style = "Normal"
$cell = $ws.Range("C37")
$cell.Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.Value = "'1.78"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = "'  -5.21%  "
$cell.Style = "Normal"
$cell = $ws.Range("B38")
$cell.Value = "'Fetch.AI"
$cell.Style = "Normal"
$cell = $ws.Range("C38")
$cell.Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.Value = "'1.23"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = "'  -3.25%  "
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.Value = "'0.131"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = "'  +3.05%  "
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = "'22.23"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "'  -0.48%  "
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "'  -0.20%  "
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "'  +0.03%  "
$cell.Style = "Normal"
$cell = $ws.Range("B43")
$cell.Value = "'Stellar"
$cell.Style = "Normal"
$cell = $ws.Range("C43")
$cell.Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = "'0.140"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'  +11.61%  "
$cell.Style = "Normal"
$cell = $ws.Range("B44")
$cell.Value = "'PolygonEcosystemToken"
$cell.Style = "Normal"
$cell = $ws.Range("C44")
$cell.Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.Value = "'0.364"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'  -0.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = "'1.81"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "'  -4.66%  "
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.Value = "'147.21"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'  +0.03%  "
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = "'0.0692"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "'  +13.61%  "
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.Value = "'43.30"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "'  -1.42%  "
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.Value = "'4.03"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "'  +1.39%  "
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = "'1.20"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "'  +0.56%  "
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.Value = "'155.14"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "'  -7.68%  "
$cell.Style = "Normal"
